# Updated code for serial run issue
# - Insert a new "CustomizeGrid" worksheet between "AddUser" and "UsersGrid"
# - Populate its header row + one data row
# - Match column widths / bold header style used elsewhere in the workbook
# - New sheet becomes the active tab (matches activeTab="1" / tabSelected="1")

$wb = $excel.ActiveWorkbook

$addUser = $wb.Worksheets.Item("AddUser")
$ws = $wb.Worksheets.Add($null, $addUser)
$ws.Name = "CustomizeGrid"

# Header row (bold, shared-string reuse of existing entries)
$ws.Cells.Item(1,1).Value = "Automation Test ID"
$ws.Cells.Item(1,2).Value = "Data"
$ws.Cells.Item(1,3).Value = "Operation"
$ws.Cells.Item(1,4).Value = "Expected Result"
$ws.Cells.Item(1,5).Value = "Actua lResult"
$ws.Cells.Item(1,6).Value = "Status"
$ws.Range("A1:F1").Font.Bold = $true

# Data row (D is written before C so new shared strings land in the same
# index order the source workbook uses: ... 37 Data, 38 Actua lResult,
# 39 Webtable customized successfully, 40 All)
$ws.Cells.Item(2,1).Value = "User_TC001"
$ws.Cells.Item(2,2).Value = "NA"
$ws.Cells.Item(2,4).Value = "Webtable customized successfully"
$ws.Cells.Item(2,3).Value = "All"

# Column widths (matches the bestFit widths used on the other sheets)
$ws.Columns.Item(1).ColumnWidth = 17.451822916666668
$ws.Columns.Item(2).ColumnWidth = 93.73697916666667
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws.Columns.Item(4).ColumnWidth = 31.166666666666668
$ws.Columns.Item(5).ColumnWidth = 11.307291666666666

$ws.Activate()
